# Insert a new data row at spreadsheet row 5, pushing existing rows 5-93
# down to 6-94 (dimension grows from A1:T93 to A1:T94), then populate the
# new row 5 with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Insert()

$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C5").Value = 'Metropolitana'
$ws.Range("D5").Value = 44882
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 'Fruta'
$ws.Range("G5").Value = 100108
$ws.Range("H5").Value = 'Tropicales y subtropicales'
$ws.Range("I5").Value = 100108007
$ws.Range("J5").Value = 'Coco'
$ws.Range("K5").Value = 'Sin especificar'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 29000
$ws.Range("Q5").Value = '$/malla 20 unidades'
$ws.Range("R5").Value = 'Perú'
$ws.Range("S5").Value = 1450
$ws.Range("T5").Value = 20
